$wb = $excel.ActiveWorkbook

# --- Sheet "progetto" (sheet1): add row 2 with project info ---
$wsProgetto = $wb.Worksheets.Item("progetto")
$wsProgetto.Range("A2").Value = "CineNow"
$wsProgetto.Range("B2").Value = "https://github.com/emanueleivn/CinemaNow"
$wsProgetto.Range("C2").Value = "WebApp per prenotazioni posti a sedere nelle sale cinema"
$wsProgetto.Columns.Item(3).ColumnWidth = 50.8984375

# --- Sheet "partecipanti" (sheet2): add rows for new team members, update project name ---
$wsPartecipanti = $wb.Worksheets.Item("partecipanti")

# Row 2: Emanuele Iovane
$wsPartecipanti.Range("A2").Value = "0512120565"
$wsPartecipanti.Range("B2").Value = "Iovane"
$wsPartecipanti.Range("C2").Value = "Emanuele"
$wsPartecipanti.Range("D2").Value = "e.iovane2@studenti.unisa.it"
$wsPartecipanti.Range("E2").Value = "CineNow"
$wsPartecipanti.Hyperlinks.Add($wsPartecipanti.Range("D2"), "mailto:e.iovane2@studenti.unisa.it") | Out-Null

# Row 3: Update project name (CinemaNow -> CineNow)
$wsPartecipanti.Range("E3").Value = "CineNow"

# Row 4: Antonio Caiazzo
$wsPartecipanti.Range("A4").Value = "0512117751"
$wsPartecipanti.Range("B4").Value = "Caiazzo"
$wsPartecipanti.Range("C4").Value = "Antonio"
$wsPartecipanti.Range("D4").Value = "a.caiazzo38@studenti.unisa.it"
$wsPartecipanti.Range("E4").Value = "CineNow"
$wsPartecipanti.Hyperlinks.Add($wsPartecipanti.Range("D4"), "mailto:a.caiazzo38@studenti.unisa.it") | Out-Null

$wsPartecipanti.Columns.Item(4).ColumnWidth = 25.59765625

# --- Active sheet / selection housekeeping ---
$wsProgetto.Range("A6").Select()
$wsProgetto.Activate()
$wsPartecipanti.Range("E4").Select()
$wsProgetto.Activate()
